$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the header labels: A1 becomes "weight", B1 becomes "bmi"
$ws.Range("A1").Value = "weight"
$ws.Range("B1").Value = "bmi"

# Update the selected cell on this sheet to B2
$ws.Range("B2").Select()
